# Scheduled runner refresh: pull latest market-board prices and re-derive
# currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) columns
# (H:N) for the leves whose underlying item prices moved since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13: The Hexster Runoff
$ws.Range("H13").Value = 2866.3333
$ws.Range("I13").Value = 2800
$ws.Range("J13").Value = 2899.5
$ws.Range("K13").Value = 2800
$ws.Range("L13").Value = 2899.5
$ws.Range("M13").Value = -2631
$ws.Range("N13").Value = -3237.5

# Row 16: Using Your Arcane Powers for Fun and Profit
$ws.Range("H16").Value = 3224.25
$ws.Range("I16").Value = 3366.3333
$ws.Range("J16").Value = 2798
$ws.Range("K16").Value = 3366.3333
$ws.Range("L16").Value = 2798
$ws.Range("M16").Value = -3136.3333
$ws.Range("N16").Value = -3258

# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 5539.1304
$ws.Range("J40").Value = 1400
$ws.Range("L40").Value = 1400
$ws.Range("N40").Value = -1750

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 3346.3242
$ws.Range("I76").Value = 3305.606
$ws.Range("J76").Value = 3682.25
$ws.Range("K76").Value = 3305.606
$ws.Range("L76").Value = 3682.25
$ws.Range("M76").Value = -2990.606
$ws.Range("N76").Value = -4312.25

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 3346.3242
$ws.Range("I79").Value = 3305.606
$ws.Range("J79").Value = 3682.25
$ws.Range("K79").Value = 3305.606
$ws.Range("L79").Value = 3682.25
$ws.Range("M79").Value = -2213.606
$ws.Range("N79").Value = -5866.25

# Row 82: Rolling on Initiative
$ws.Range("H82").Value = 3428.4348
$ws.Range("I82").Value = 1179.3125
$ws.Range("J82").Value = 8569.286
$ws.Range("K82").Value = 3537.9375
$ws.Range("L82").Value = 25707.858
$ws.Range("M82").Value = -3131.9375
$ws.Range("N82").Value = -26519.858

# Row 85: Darkly Dreaming Dexterity (L)
$ws.Range("H85").Value = 3428.4348
$ws.Range("I85").Value = 1179.3125
$ws.Range("J85").Value = 8569.286
$ws.Range("K85").Value = 3537.9375
$ws.Range("L85").Value = 25707.858
$ws.Range("M85").Value = -2133.9375
$ws.Range("N85").Value = -28515.858

# Row 98: The Dotted Line
$ws.Range("H98").Value = 2125.45
$ws.Range("I98").Value = 1220.1875
$ws.Range("J98").Value = 5746.5
$ws.Range("K98").Value = 1220.1875
$ws.Range("L98").Value = 5746.5
$ws.Range("M98").Value = 277.8125
$ws.Range("N98").Value = -8742.5

# Row 106: Making Your Mark
$ws.Range("H106").Value = 3044.4707
$ws.Range("J106").Value = 3016
$ws.Range("L106").Value = 3016
$ws.Range("N106").Value = -4278

# Row 109: A Time for Peace
$ws.Range("H109").Value = 33684
$ws.Range("J109").Value = 33684
$ws.Range("L109").Value = 33684
$ws.Range("N109").Value = -36458

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 975
$ws.Range("I111").Value = 975
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2925
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 142
$ws.Range("N111").ClearContents()

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 1539.7954
$ws.Range("J112").Value = 1780.6666
$ws.Range("L112").Value = 5341.9998
$ws.Range("N112").Value = -7557.9998

# Row 115: 5-bell Energy
$ws.Range("H115").Value = 1871.6666
$ws.Range("I115").Value = 246
$ws.Range("J115").Value = 10000
$ws.Range("K115").Value = 738
$ws.Range("L115").Value = 30000
$ws.Range("M115").Value = 829
$ws.Range("N115").Value = -33134

# Row 118: Crafty Concoctions
$ws.Range("H118").Value = 791.6667
$ws.Range("I118").Value = 350
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 1050
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = 607
$ws.Range("N118").Value = -12314

# Row 122: Wishful Inking
$ws.Range("H122").Value = 2125.45
$ws.Range("I122").Value = 1220.1875
$ws.Range("J122").Value = 5746.5
$ws.Range("K122").Value = 3660.5625
$ws.Range("L122").Value = 17239.5
$ws.Range("M122").Value = -1210.5625
$ws.Range("N122").Value = -22139.5

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 956.7692
$ws.Range("I74").Value = 787.75
$ws.Range("J74").Value = 2985
$ws.Range("K74").Value = 787.75
$ws.Range("L74").Value = 2985
$ws.Range("M74").Value = 86.25
$ws.Range("N74").Value = -4733

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 956.7692
$ws.Range("I77").Value = 787.75
$ws.Range("J77").Value = 2985
$ws.Range("K77").Value = 3938.75
$ws.Range("L77").Value = 14925
$ws.Range("M77").Value = 429.25
$ws.Range("N77").Value = -23661

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 5195.4443
$ws.Range("I102").Value = 5251.2856
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 5251.2856
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -3629.2856
$ws.Range("N102").Value = -8244

$ws = $wb.Worksheets.Item("BSM")
# Row 87: Winter Weather Conditions
$ws.Range("H87").Value = 45344.08
$ws.Range("J87").Value = 45344.08
$ws.Range("L87").Value = 45344.08
$ws.Range("N87").Value = -47840.08

# Row 90: The Nightsoil Is Dark and Full of Terrors (L)
$ws.Range("H90").Value = 45344.08
$ws.Range("J90").Value = 45344.08
$ws.Range("L90").Value = 136032.24
$ws.Range("N90").Value = -148512.24

# Row 107: The Gold Experience
$ws.Range("H107").Value = 2181.6667
$ws.Range("I107").Value = 2203.6667
$ws.Range("J107").Value = 2137.6667
$ws.Range("K107").Value = 2203.6667
$ws.Range("L107").Value = 2137.6667
$ws.Range("M107").Value = -283.6667000000002
$ws.Range("N107").Value = -5977.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 5801.2104
$ws.Range("I31").Value = 10764.8
$ws.Range("J31").Value = 4028.5
$ws.Range("K31").Value = 10764.8
$ws.Range("L31").Value = 4028.5
$ws.Range("M31").Value = -10469.8
$ws.Range("N31").Value = -4618.5

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 5801.2104
$ws.Range("I34").Value = 10764.8
$ws.Range("J34").Value = 4028.5
$ws.Range("K34").Value = 10764.8
$ws.Range("L34").Value = 4028.5
$ws.Range("M34").Value = -10562.8
$ws.Range("N34").Value = -4432.5

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 4129.294
$ws.Range("I122").Value = 4295.263
$ws.Range("J122").Value = 3919.0667
$ws.Range("K122").Value = 12885.789
$ws.Range("L122").Value = 11757.2001
$ws.Range("M122").Value = -10435.789
$ws.Range("N122").Value = -16657.2001

$ws = $wb.Worksheets.Item("CUL")
# Row 110: His Dark Utensils
$ws.Range("H110").Value = 3468.375
$ws.Range("I110").Value = 1949.4
$ws.Range("K110").Value = 5848.200000000001
$ws.Range("M110").Value = -1758.200000000001

# Row 125: At Any Temperature
$ws.Range("H125").Value = 6388.8335
$ws.Range("I125").Value = 4250
$ws.Range("J125").Value = 7458.25
$ws.Range("K125").Value = 12750
$ws.Range("L125").Value = 22374.75
$ws.Range("M125").Value = -7830
$ws.Range("N125").Value = -32214.75

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 855.91895
$ws.Range("I131").Value = 447.27274
$ws.Range("J131").Value = 1028.8077
$ws.Range("K131").Value = 1341.81822
$ws.Range("L131").Value = 3086.4231
$ws.Range("M131").Value = 3698.18178
$ws.Range("N131").Value = -13166.4231

$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 1292.8
$ws.Range("I126").Value = 1154.6666
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 3463.9998
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -993.9998000000001
$ws.Range("N126").Value = -9440

# Row 132: On Board for Lar
$ws.Range("H132").Value = 3170.1333
$ws.Range("I132").Value = 3188
$ws.Range("J132").Value = 3161.2
$ws.Range("K132").Value = 9564
$ws.Range("L132").Value = 9483.599999999999
$ws.Range("M132").Value = -7034
$ws.Range("N132").Value = -14543.6

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 2327.5881
$ws.Range("I40").Value = 2145.8518
$ws.Range("J40").Value = 3028.5715
$ws.Range("K40").Value = 2145.8518
$ws.Range("L40").Value = 3028.5715
$ws.Range("M40").Value = -2009.8518
$ws.Range("N40").Value = -3300.5715

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 528158.3
$ws.Range("I122").Value = 1112467
$ws.Range("K122").Value = 3337401
$ws.Range("M122").Value = -3334951

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 435887.9
$ws.Range("I126").Value = 556350.9399999999
$ws.Range("J126").Value = 2221
$ws.Range("K126").Value = 1669052.82
$ws.Range("L126").Value = 6663
$ws.Range("M126").Value = -1666582.82
$ws.Range("N126").Value = -11603
